$wb = $excel.ActiveWorkbook

# --- Sheet: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.9395017793594306
$ws1.Range("C2").Value = 0.4
$ws1.Range("D2").Value = 0.4285714285714285
$ws1.Range("E2").Value = 0.4137931034482759
$ws1.Range("F2").Value = 0.4225352112676056
$ws1.Range("G2").Value = 0.4273972602739726
$ws1.Range("H2").Value = 0.6974317817014447
$ws1.Range("I2").Value = 12
$ws1.Range("J2").Value = 18
$ws1.Range("K2").Value = 516
$ws1.Range("L2").Value = 16

# --- Sheet: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")

# Row 2 (label "0")
$ws2.Range("B2").Value = 0.9699248120300752
$ws2.Range("C2").Value = 0.9662921348314607
$ws2.Range("D2").Value = 0.9681050656660413

# Row 3 (label "1")
$ws2.Range("B3").Value = 0.4
$ws2.Range("C3").Value = 0.4285714285714285
$ws2.Range("D3").Value = 0.4137931034482759

# Row 4 (label "accuracy")
$ws2.Range("B4").Value = 0.9395017793594306
$ws2.Range("C4").Value = 0.9395017793594306
$ws2.Range("D4").Value = 0.9395017793594306
$ws2.Range("E4").Value = 0.9395017793594306

# Row 5 (label "macro avg")
$ws2.Range("B5").Value = 0.6849624060150377
$ws2.Range("C5").Value = 0.6974317817014446
$ws2.Range("D5").Value = 0.6909490845571585

# Row 6 (label "weighted avg")
$ws2.Range("B6").Value = 0.9415299815374736
$ws2.Range("C6").Value = 0.9395017793594306
$ws2.Range("D6").Value = 0.9404880995768999

# --- Sheet: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

$ws3.Range("B2").Value = 516
$ws3.Range("C2").Value = 18

$ws3.Range("B3").Value = 16
$ws3.Range("C3").Value = 12
